$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: move the "x" mark from C5 to D5 (label "Inicio de sesión" unchanged)
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "x"

# Row 9: move the "x" mark from F9 to C9 (label "Vender" unchanged)
$ws.Range("F9").Value = $null
$ws.Range("C9").Value = "x"

# Row 11: move the "x" mark from E11 to D11 (label "Página de los productos" unchanged)
$ws.Range("E11").Value = $null
$ws.Range("D11").Value = "x"

# Row 13: relabel "Seccion de anime" -> "Misproductos" and move "x" from F13 to C13
$ws.Range("A13").Value = "Misproductos"
$ws.Range("F13").Value = $null
$ws.Range("C13").Value = "x"

# Row 14: relabel "Favoritos" -> "Filtro por secciones" ("x" stays in F14)
$ws.Range("A14").Value = "Filtro por secciones"

# Row 15: clear out the row entirely (label "Seccion de electrónica" and "x" removed)
$ws.Range("A15").Value = $null
$ws.Range("F15").Value = $null

# Update the active selection to A15:B15
$ws.Range("A15:B15").Select()
